$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = 1
$ws.Range("N1").Value = "2"
$ws.Range("O1").Value = "3"
$ws.Range("P1").Value = "4"
$ws.Range("Q1").Value = "5"
$ws.Range("R1").Value = "6"
$ws.Range("S1").Value = "7"
$ws.Range("T1").Value = "8"
$ws.Range("U1").Value = "9"
$ws.Range("V1").Value = "10"
